$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 86, shifting existing rows 86-175 down to 89-178
$ws.Rows.Item(86).Resize(3).Insert()

# Fill the 3 new rows (86, 87, 88) with the new week's data.
$rows = @(
    @{ Row=86; D=44512; L="Primera"; M=250; N=9000;  O=9000;  P=9000; Q="`$/bandeja 7 kilos"; R="Provincia de Melipilla";   S=1286 },
    @{ Row=87; D=44512; L="Primera"; M=185; N=9000;  O=10000; P=9595; Q="`$/bandeja 7 kilos"; R="Región de La Araucanía";  S=1371 },
    @{ Row=88; D=44512; L="Segunda"; M=150; N=7000;  O=7000;  P=7000; Q="`$/bandeja 7 kilos"; R="Provincia de Melipilla";   S=1000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 10
    $ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
    $ws.Cells.Item($row, 3).Value = "La Araucanía"
    $ws.Cells.Item($row, 4).Value = [DateTime]::FromOADate($r.D)
    $ws.Cells.Item($row, 5).Value = 9
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 7
}
